$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains its original text representation
# (values like "29.764.54" or "1.000" must not be reinterpreted as numbers)
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '29.764.54'
$ws.Range('E2').Value = '  +2.11%  '

# Row 3
$ws.Range('D3').Value = '1.857.31'
$ws.Range('E3').Value = '  +1.76%  '

# Row 4
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  +0.08%  '

# Row 5
$ws.Range('D5').Value = '245.36'
$ws.Range('E5').Value = '  +1.56%  '

# Row 6
$ws.Range('D6').Value = '0.6427'
$ws.Range('E6').Value = '  +3.50%  '

# Row 7
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').Value = '46.91'
$ws.Range('E8').Value = '  +4.07%  '

# Row 9
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = '0.07499'
$ws.Range('E9').Value = '  +2.24%  '

# Row 10
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').Value = '0.2991'
$ws.Range('E10').Value = '  +3.08%  '

# Row 11
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D11').Value = '24.08'
$ws.Range('E11').Value = '  +4.71%  '

# Row 12
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = '0.07688'
$ws.Range('E12').Value = '  +0.12%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.872.42'
$ws.Range('E13').Value = '  +2.68%  '

# Row 14
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '5.051'
$ws.Range('E14').Value = '  +1.88%  '

# Row 15
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = '0.6838'
$ws.Range('E15').Value = '  +2.88%  '

# Row 16
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '83.93'
$ws.Range('E16').Value = '  +1.97%  '

# Row 17
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.000009431'
$ws.Range('E17').Value = '  +5.11%  '

# Row 18
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = '6.062'
$ws.Range('E18').Value = '  +3.73%  '

# Row 19
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '29.747.59'
$ws.Range('E19').Value = '  +2.14%  '

# Row 20
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').Value = '2.122.75'
$ws.Range('E20').Value = '  +2.60%  '

# Row 21
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '240.04'
$ws.Range('E21').Value = '  +0.62%  '

# Row 22
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').Value = '12.68'
$ws.Range('E22').Value = '  +2.14%  '

# Row 23
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  +0.06%  '

# Row 24
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').Value = '7.452'
$ws.Range('E24').Value = '  +2.44%  '

# Row 25
$ws.Range('B25').Value = 'BinanceUSD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D25').Value = '1.001'
$ws.Range('E25').Value = '  +0.10%  '

# Row 26
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '158.89'
$ws.Range('E26').Value = '  +0.52%  '

# Row 27
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = '0.1428'
$ws.Range('E27').Value = '  +0.23%  '

# Row 28
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '8.528'
$ws.Range('E28').Value = '  +0.55%  '

# Row 29
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '17.96'
$ws.Range('E29').Value = '  +1.73%  '

# Row 30
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = '0.06168'
$ws.Range('E30').Value = '  +10.80%  '

# Row 31
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.504'
$ws.Range('E31').Value = '  +1.26%  '

# Row 32
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = '1.268'
$ws.Range('E32').Value = '  +5.29%  '

# Row 33
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '4.150'
$ws.Range('E33').Value = '  +1.39%  '

# Row 34
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = '4.124'
$ws.Range('E34').Value = '  +0.92%  '

# Row 35
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = '1.882'
$ws.Range('E35').Value = '  +2.32%  '

# Row 36
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '1.158'
$ws.Range('E36').Value = '  +2.13%  '

# Row 37
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '0.7334'
$ws.Range('E37').Value = '  -0.01%  '

# Row 38
$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').Value = '2.606'
$ws.Range('E38').Value = '  -0.73%  '

# Row 39
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '2.852'
$ws.Range('E39').Value = '  +0.35%  '

# Row 40
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '0.01792'
$ws.Range('E40').Value = '  +1.56%  '

# Row 41
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '1.218.82'
$ws.Range('E41').Value = '  +0.62%  '

# Row 42
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '0.9292'
$ws.Range('E42').Value = '  +1.81%  '

# Row 43
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '6.205'
$ws.Range('E43').Value = '  -1.58%  '

# Row 44
$ws.Range('D44').Value = '2.034.40'
$ws.Range('E44').Value = '  +3.16%  '

# Row 45
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  +0.11%  '

# Row 46
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = '101.88'
$ws.Range('E46').Value = '  +0.21%  '

# Row 47
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '66.47'
$ws.Range('E47').Value = '  +2.92%  '

# Row 48
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.00000000121'
$ws.Range('E48').Value = '  +3.19%  '

# Row 49
$ws.Range('D49').Value = '9.344'
$ws.Range('E49').Value = '  +2.29%  '

# Row 50
$ws.Range('D50').Value = '0.4081'
$ws.Range('E50').Value = '  +1.45%  '

# Row 51
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.5007'
$ws.Range('E51').Value = '  -1.56%  '
